$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Save the current B2 value; it moves down to become B3's value.
$oldB2 = $ws.Range("B2").Value2

# Copy A2's formatting (style) to the new A3 cell, then set its value.
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A3").Value = 1

# B3 receives the value that used to be in B2.
$ws.Range("B3").Value = $oldB2

# B2 gets a new value.
$ws.Range("B2").Value = 33.54101966249685
